# Auto-generated Excel COM-interop script
# Applies cached-value updates (market price refresh) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 1138.4117  # H18: 1086.2222 -> 1138.4117
$ws.Cells.Item(18, 9).Value = 1054.0834  # I18: 988.3077 -> 1054.0834
$ws.Cells.Item(18, 11).Value = 1054.0834  # K18: 988.3077 -> 1054.0834
$ws.Cells.Item(18, 13).Value = -770.0834  # M18: -704.3077 -> -770.0834
$ws.Cells.Item(40, 8).Value = 2099.8  # H40: 2049.8333 -> 2099.8
$ws.Cells.Item(40, 9).Value = 2233  # I40: 2124.75 -> 2233
$ws.Cells.Item(40, 11).Value = 2233  # K40: 2124.75 -> 2233
$ws.Cells.Item(40, 13).Value = -2058  # M40: -1949.75 -> -2058
$ws.Cells.Item(43, 8).Value = 6945211  # H43: 13889989 -> 6945211
$ws.Cells.Item(43, 9).Value = 1033.5  # I43: 1700.5 -> 1033.5
$ws.Cells.Item(43, 10).Value = 13889389  # J43: 27778278 -> 13889389
$ws.Cells.Item(43, 11).Value = 1033.5  # K43: 1700.5 -> 1033.5
$ws.Cells.Item(43, 12).Value = 13889389  # L43: 27778278 -> 13889389
$ws.Cells.Item(43, 13).Value = -964.5  # M43: -1631.5 -> -964.5
$ws.Cells.Item(43, 14).Value = -13889527  # N43: -27778416 -> -13889527
$ws.Cells.Item(51, 8).Value = 2676  # H51: 2678 -> 2676
$ws.Cells.Item(51, 10).Value = 2895  # J51: 2897.5 -> 2895
$ws.Cells.Item(51, 12).Value = 2895  # L51: 2897.5 -> 2895
$ws.Cells.Item(51, 14).Value = -3863  # N51: -3865.5 -> -3863
$ws.Cells.Item(61, 8).Value = 50  # H61: 47.5 -> 50
$ws.Cells.Item(61, 9).Value = 50  # I61: 47.5 -> 50
$ws.Cells.Item(61, 11).Value = 150  # K61: 142.5 -> 150
$ws.Cells.Item(61, 13).Value = 22  # M61: 29.5 -> 22
$ws.Cells.Item(100, 8).Value = 1904.762  # H100: 1604.1111 -> 1904.762
$ws.Cells.Item(100, 9).Value = 1700  # I100: 1184.4445 -> 1700
$ws.Cells.Item(100, 10).Value = 1968.75  # J100: 1813.9445 -> 1968.75
$ws.Cells.Item(100, 11).Value = 1700  # K100: 1184.4445 -> 1700
$ws.Cells.Item(100, 12).Value = 1968.75  # L100: 1813.9445 -> 1968.75
$ws.Cells.Item(100, 13).Value = -1159  # M100: -643.4445000000001 -> -1159
$ws.Cells.Item(100, 14).Value = -3050.75  # N100: -2895.9445 -> -3050.75
$ws.Cells.Item(111, 8).Value = 1664.8  # H111: 1626.8572 -> 1664.8
$ws.Cells.Item(111, 9).Value = 1581  # I111: 1641.3334 -> 1581
$ws.Cells.Item(111, 10).Value = 2000  # J111: 1616 -> 2000
$ws.Cells.Item(111, 11).Value = 4743  # K111: 4924.0002 -> 4743
$ws.Cells.Item(111, 12).Value = 6000  # L111: 4848 -> 6000
$ws.Cells.Item(111, 13).Value = -1676  # M111: -1857.0002 -> -1676
$ws.Cells.Item(111, 14).Value = -12134  # N111: -10982 -> -12134
$ws.Cells.Item(132, 8).Value = 8338446.5  # H132: 9529570 -> 8338446.5
$ws.Cells.Item(132, 9).Value = 11117065  # I132: 12827294 -> 11117065
$ws.Cells.Item(132, 10).Value = 2590.3  # J132: 2810.7778 -> 2590.3
$ws.Cells.Item(132, 11).Value = 33351195  # K132: 38481882 -> 33351195
$ws.Cells.Item(132, 12).Value = 7770.900000000001  # L132: 8432.3334 -> 7770.900000000001
$ws.Cells.Item(132, 13).Value = -33348665  # M132: -38479352 -> -33348665
$ws.Cells.Item(132, 14).Value = -12830.9  # N132: -13492.3334 -> -12830.9
$ws.Cells.Item(137, 8).Value = 1064.7556  # H137: 1055.775 -> 1064.7556
$ws.Cells.Item(137, 9).Value = 684.96295  # I137: 671.8214 -> 684.96295
$ws.Cells.Item(137, 10).Value = 1634.4445  # J137: 1951.6666 -> 1634.4445
$ws.Cells.Item(137, 11).Value = 2054.88885  # K137: 2015.4642 -> 2054.88885
$ws.Cells.Item(137, 12).Value = 4903.333500000001  # L137: 5854.9998 -> 4903.333500000001
$ws.Cells.Item(137, 13).Value = 495.1111500000002  # M137: 534.5357999999999 -> 495.1111500000002
$ws.Cells.Item(137, 14).Value = -10003.3335  # N137: -10954.9998 -> -10003.3335
$ws.Cells.Item(141, 8).Value = 555.7931  # H141: 568.7692 -> 555.7931
$ws.Cells.Item(141, 9).Value = 555.7931  # I141: 568.7692 -> 555.7931
$ws.Cells.Item(141, 11).Value = 1667.3793  # K141: 1706.3076 -> 1667.3793
$ws.Cells.Item(141, 13).Value = 3512.6207  # M141: 3473.6924 -> 3512.6207

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5608.9  # H2: 6887.9375 -> 5608.9
$ws.Cells.Item(2, 9).Value = 583.2857  # I2: 619.5 -> 583.2857
$ws.Cells.Item(2, 11).Value = 583.2857  # K2: 619.5 -> 583.2857
$ws.Cells.Item(2, 13).Value = -470.2857  # M2: -506.5 -> -470.2857
$ws.Cells.Item(25, 8).Value = 1700  # H25: 6500 -> 1700
$ws.Cells.Item(25, 9).Value = 2000  # I25: 3000 -> 2000
$ws.Cells.Item(25, 10).Value = 1400  # J25: 10000 -> 1400
$ws.Cells.Item(25, 11).Value = 2000  # K25: 3000 -> 2000
$ws.Cells.Item(25, 12).Value = 1400  # L25: 10000 -> 1400
$ws.Cells.Item(25, 13).Value = -1598  # M25: -2598 -> -1598
$ws.Cells.Item(25, 14).Value = -2204  # N25: -10804 -> -2204
$ws.Cells.Item(45, 8).Value = 1236  # H45: 1082.05 -> 1236
$ws.Cells.Item(45, 9).Value = 1270.2858  # I45: 954 -> 1270.2858
$ws.Cells.Item(45, 10).Value = 1206  # J45: 1238.5555 -> 1206
$ws.Cells.Item(45, 11).Value = 1270.2858  # K45: 954 -> 1270.2858
$ws.Cells.Item(45, 12).Value = 1206  # L45: 1238.5555 -> 1206
$ws.Cells.Item(45, 13).Value = -893.2858000000001  # M45: -577 -> -893.2858000000001
$ws.Cells.Item(45, 14).Value = -1960  # N45: -1992.5555 -> -1960
$ws.Cells.Item(110, 8).Value = 1804.2858  # H110: 2269.4285 -> 1804.2858
$ws.Cells.Item(110, 9).Value = 1480.9166  # I110: 1869.875 -> 1480.9166
$ws.Cells.Item(110, 10).Value = 2235.4443  # J110: 2802.1667 -> 2235.4443
$ws.Cells.Item(110, 11).Value = 1480.9166  # K110: 1869.875 -> 1480.9166
$ws.Cells.Item(110, 12).Value = 2235.4443  # L110: 2802.1667 -> 2235.4443
$ws.Cells.Item(110, 13).Value = 564.0834  # M110: 175.125 -> 564.0834
$ws.Cells.Item(110, 14).Value = -6325.4443  # N110: -6892.1667 -> -6325.4443
$ws.Cells.Item(116, 8).Value = 5608.9  # H116: 6887.9375 -> 5608.9
$ws.Cells.Item(116, 9).Value = 583.2857  # I116: 619.5 -> 583.2857
$ws.Cells.Item(116, 11).Value = 583.2857  # K116: 619.5 -> 583.2857
$ws.Cells.Item(116, 13).Value = 1710.7143  # M116: 1674.5 -> 1710.7143
$ws.Cells.Item(132, 8).Value = 2056.75  # H132: 2131.1924 -> 2056.75
$ws.Cells.Item(132, 9).Value = 1817  # I132: 1892.8 -> 1817
$ws.Cells.Item(132, 10).Value = 2376.4167  # J132: 2456.2727 -> 2376.4167
$ws.Cells.Item(132, 11).Value = 5451  # K132: 5678.4 -> 5451
$ws.Cells.Item(132, 12).Value = 7129.250100000001  # L132: 7368.8181 -> 7129.250100000001
$ws.Cells.Item(132, 13).Value = -2921  # M132: -3148.4 -> -2921
$ws.Cells.Item(132, 14).Value = -12189.2501  # N132: -12428.8181 -> -12189.2501

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5608.9  # H3: 6887.9375 -> 5608.9
$ws.Cells.Item(3, 9).Value = 583.2857  # I3: 619.5 -> 583.2857
$ws.Cells.Item(3, 11).Value = 583.2857  # K3: 619.5 -> 583.2857
$ws.Cells.Item(3, 13).Value = -469.2857  # M3: -505.5 -> -469.2857
$ws.Cells.Item(80, 8).Value = 919.625  # H80: 935.25 -> 919.625
$ws.Cells.Item(80, 9).Value = 727.4286  # I80: 742.75 -> 727.4286
$ws.Cells.Item(80, 10).Value = 1069.1111  # J80: 1127.75 -> 1069.1111
$ws.Cells.Item(80, 11).Value = 727.4286  # K80: 742.75 -> 727.4286
$ws.Cells.Item(80, 12).Value = 1069.1111  # L80: 1127.75 -> 1069.1111
$ws.Cells.Item(80, 13).Value = 270.5714  # M80: 255.25 -> 270.5714
$ws.Cells.Item(80, 14).Value = -3065.1111  # N80: -3123.75 -> -3065.1111
$ws.Cells.Item(83, 8).Value = 919.625  # H83: 935.25 -> 919.625
$ws.Cells.Item(83, 9).Value = 727.4286  # I83: 742.75 -> 727.4286
$ws.Cells.Item(83, 10).Value = 1069.1111  # J83: 1127.75 -> 1069.1111
$ws.Cells.Item(83, 11).Value = 3637.143  # K83: 3713.75 -> 3637.143
$ws.Cells.Item(83, 12).Value = 5345.5555  # L83: 5638.75 -> 5345.5555
$ws.Cells.Item(83, 13).Value = 1354.857  # M83: 1278.25 -> 1354.857
$ws.Cells.Item(83, 14).Value = -15329.5555  # N83: -15622.75 -> -15329.5555
$ws.Cells.Item(134, 8).Value = 8493.177  # H134: 10032.143 -> 8493.177
$ws.Cells.Item(134, 9).Value = 1647.1  # I134: 1731.375 -> 1647.1
$ws.Cells.Item(134, 10).Value = 18273.285  # J134: 21099.834 -> 18273.285
$ws.Cells.Item(134, 11).Value = 4941.299999999999  # K134: 5194.125 -> 4941.299999999999
$ws.Cells.Item(134, 12).Value = 54819.855  # L134: 63299.50199999999 -> 54819.855
$ws.Cells.Item(134, 13).Value = -2406.299999999999  # M134: -2659.125 -> -2406.299999999999
$ws.Cells.Item(134, 14).Value = -59889.855  # N134: -68369.50199999999 -> -59889.855

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1474.3846  # H31: 1719.2 -> 1474.3846
$ws.Cells.Item(31, 9).Value = 1108.6  # I31: 1345.6 -> 1108.6
$ws.Cells.Item(31, 10).Value = 1973.1818  # J31: 2092.8 -> 1973.1818
$ws.Cells.Item(31, 11).Value = 1108.6  # K31: 1345.6 -> 1108.6
$ws.Cells.Item(31, 12).Value = 1973.1818  # L31: 2092.8 -> 1973.1818
$ws.Cells.Item(31, 13).Value = -813.5999999999999  # M31: -1050.6 -> -813.5999999999999
$ws.Cells.Item(31, 14).Value = -2563.1818  # N31: -2682.8 -> -2563.1818
$ws.Cells.Item(34, 8).Value = 1474.3846  # H34: 1719.2 -> 1474.3846
$ws.Cells.Item(34, 9).Value = 1108.6  # I34: 1345.6 -> 1108.6
$ws.Cells.Item(34, 10).Value = 1973.1818  # J34: 2092.8 -> 1973.1818
$ws.Cells.Item(34, 11).Value = 1108.6  # K34: 1345.6 -> 1108.6
$ws.Cells.Item(34, 12).Value = 1973.1818  # L34: 2092.8 -> 1973.1818
$ws.Cells.Item(34, 13).Value = -906.5999999999999  # M34: -1143.6 -> -906.5999999999999
$ws.Cells.Item(34, 14).Value = -2377.1818  # N34: -2496.8 -> -2377.1818
$ws.Cells.Item(58, 8).Value = 1134.5  # H58: 1012.1818 -> 1134.5
$ws.Cells.Item(58, 9).Value = 1096.8  # I58: 979.4286 -> 1096.8
$ws.Cells.Item(58, 11).Value = 1096.8  # K58: 979.4286 -> 1096.8
$ws.Cells.Item(58, 13).Value = -893.8  # M58: -776.4286 -> -893.8
$ws.Cells.Item(86, 8).Value = 6082499.5  # H86: 7433276 -> 6082499.5
$ws.Cells.Item(86, 9).Value = 22225090  # I86: 33335636 -> 22225090
$ws.Cells.Item(86, 10).Value = 29027.75  # J86: 32602 -> 29027.75
$ws.Cells.Item(86, 11).Value = 22225090  # K86: 33335636 -> 22225090
$ws.Cells.Item(86, 12).Value = 29027.75  # L86: 32602 -> 29027.75
$ws.Cells.Item(86, 13).Value = -22223967  # M86: -33334513 -> -22223967
$ws.Cells.Item(86, 14).Value = -31273.75  # N86: -34848 -> -31273.75
$ws.Cells.Item(89, 8).Value = 6082499.5  # H89: 7433276 -> 6082499.5
$ws.Cells.Item(89, 9).Value = 22225090  # I89: 33335636 -> 22225090
$ws.Cells.Item(89, 10).Value = 29027.75  # J89: 32602 -> 29027.75
$ws.Cells.Item(89, 11).Value = 111125450  # K89: 166678180 -> 111125450
$ws.Cells.Item(89, 12).Value = 145138.75  # L89: 163010 -> 145138.75
$ws.Cells.Item(89, 13).Value = -111119834  # M89: -166672564 -> -111119834
$ws.Cells.Item(89, 14).Value = -156370.75  # N89: -174242 -> -156370.75
$ws.Cells.Item(132, 8).Value = 6400.913  # H132: 7288.45 -> 6400.913
$ws.Cells.Item(132, 9).Value = 7523.3125  # I132: 9826.75 -> 7523.3125
$ws.Cells.Item(132, 10).Value = 3835.4285  # J132: 3481 -> 3835.4285
$ws.Cells.Item(132, 11).Value = 22569.9375  # K132: 29480.25 -> 22569.9375
$ws.Cells.Item(132, 12).Value = 11506.2855  # L132: 10443 -> 11506.2855
$ws.Cells.Item(132, 13).Value = -20039.9375  # M132: -26950.25 -> -20039.9375
$ws.Cells.Item(132, 14).Value = -16566.2855  # N132: -15503 -> -16566.2855
$ws.Cells.Item(134, 8).Value = 2069.1052  # H134: 1959.35 -> 2069.1052
$ws.Cells.Item(134, 9).Value = 2275.1428  # I134: 2035.8125 -> 2275.1428
$ws.Cells.Item(134, 10).Value = 1492.2  # J134: 1653.5 -> 1492.2
$ws.Cells.Item(134, 11).Value = 6825.428400000001  # K134: 6107.4375 -> 6825.428400000001
$ws.Cells.Item(134, 12).Value = 4476.6  # L134: 4960.5 -> 4476.6
$ws.Cells.Item(134, 13).Value = -4290.428400000001  # M134: -3572.4375 -> -4290.428400000001
$ws.Cells.Item(134, 14).Value = -9546.6  # N134: -10030.5 -> -9546.6
$ws.Cells.Item(136, 8).Value = 1134.5  # H136: 1012.1818 -> 1134.5
$ws.Cells.Item(136, 9).Value = 1096.8  # I136: 979.4286 -> 1096.8
$ws.Cells.Item(136, 11).Value = 3290.4  # K136: 2938.2858 -> 3290.4
$ws.Cells.Item(136, 13).Value = -740.3999999999996  # M136: -388.2857999999997 -> -740.3999999999996

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1243.0312  # H5: 1274.0968 -> 1243.0312
$ws.Cells.Item(5, 9).Value = 1299.1786  # I5: 1336.9259 -> 1299.1786
$ws.Cells.Item(5, 11).Value = 3897.5358  # K5: 4010.7777 -> 3897.5358
$ws.Cells.Item(5, 13).Value = -3785.5358  # M5: -3898.7777 -> -3785.5358
$ws.Cells.Item(107, 8).Value = 9604.182000000001  # H107: 8078.6924 -> 9604.182000000001
$ws.Cells.Item(107, 9).Value = 625  # I107: 296.8 -> 625
$ws.Cells.Item(107, 10).Value = 11599.556  # J107: 12942.375 -> 11599.556
$ws.Cells.Item(107, 11).Value = 1875  # K107: 890.4000000000001 -> 1875
$ws.Cells.Item(107, 12).Value = 34798.66800000001  # L107: 38827.125 -> 34798.66800000001
$ws.Cells.Item(107, 13).Value = 45  # M107: 1029.6 -> 45
$ws.Cells.Item(107, 14).Value = -38638.66800000001  # N107: -42667.125 -> -38638.66800000001
$ws.Cells.Item(131, 8).Value = 11630119  # H131: 1223.9691 -> 11630119
$ws.Cells.Item(131, 9).Value = 1000000000  # I131: 499.5 -> 1000000000
$ws.Cells.Item(131, 10).Value = 2237.9883  # J131: 1239.2211 -> 2237.9883
$ws.Cells.Item(131, 11).Value = 3000000000  # K131: 1498.5 -> 3000000000
$ws.Cells.Item(131, 12).Value = 6713.9649  # L131: 3717.6633 -> 6713.9649
$ws.Cells.Item(131, 13).Value = -2999994960  # M131: 3541.5 -> -2999994960
$ws.Cells.Item(131, 14).Value = -16793.9649  # N131: -13797.6633 -> -16793.9649
$ws.Cells.Item(132, 8).Value = 1833.6666  # H132: 1773.091 -> 1833.6666
$ws.Cells.Item(132, 10).Value = 3833.3333  # J132: 4500 -> 3833.3333
$ws.Cells.Item(132, 12).Value = 34499.9997  # L132: 40500 -> 34499.9997
$ws.Cells.Item(132, 14).Value = -39559.9997  # N132: -45560 -> -39559.9997
$ws.Cells.Item(135, 8).Value = 1243.0312  # H135: 1274.0968 -> 1243.0312
$ws.Cells.Item(135, 9).Value = 1299.1786  # I135: 1336.9259 -> 1299.1786
$ws.Cells.Item(135, 11).Value = 11692.6074  # K135: 12032.3331 -> 11692.6074
$ws.Cells.Item(135, 13).Value = -9157.607399999999  # M135: -9497.3331 -> -9157.607399999999
$ws.Cells.Item(136, 8).Value = 1687.1765  # H136: 1747.0714 -> 1687.1765
$ws.Cells.Item(136, 9).Value = 1157.6923  # I136: 1260 -> 1157.6923
$ws.Cells.Item(136, 10).Value = 3408  # J136: 3533 -> 3408
$ws.Cells.Item(136, 11).Value = 3473.0769  # K136: 3780 -> 3473.0769
$ws.Cells.Item(136, 12).Value = 10224  # L136: 10599 -> 10224
$ws.Cells.Item(136, 13).Value = 1626.9231  # M136: 1320 -> 1626.9231
$ws.Cells.Item(136, 14).Value = -20424  # N136: -20799 -> -20424
$ws.Cells.Item(141, 8).Value = 3332.5  # H141: 3666.6667 -> 3332.5
$ws.Cells.Item(141, 9).Value = 2943.3333  # I141: 3250 -> 2943.3333
$ws.Cells.Item(141, 11).Value = 8829.999899999999  # K141: 9750 -> 8829.999899999999
$ws.Cells.Item(141, 13).Value = -3649.999899999999  # M141: -4570 -> -3649.999899999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1816.8636  # H122: 1993.7894 -> 1816.8636
$ws.Cells.Item(122, 9).Value = 1645.8422  # I122: 1845.4667 -> 1645.8422
$ws.Cells.Item(122, 10).Value = 2900  # J122: 2550 -> 2900
$ws.Cells.Item(122, 11).Value = 4937.5266  # K122: 5536.4001 -> 4937.5266
$ws.Cells.Item(122, 12).Value = 8700  # L122: 7650 -> 8700
$ws.Cells.Item(122, 13).Value = -2487.5266  # M122: -3086.4001 -> -2487.5266
$ws.Cells.Item(122, 14).Value = -13600  # N122: -12550 -> -13600
$ws.Cells.Item(126, 8).Value = 2188.889  # H126: 2113.6365 -> 2188.889
$ws.Cells.Item(126, 9).Value = 1783.3334  # I126: 1778.5714 -> 1783.3334
$ws.Cells.Item(126, 10).Value = 3000  # J126: 2700 -> 3000
$ws.Cells.Item(126, 11).Value = 5350.0002  # K126: 5335.7142 -> 5350.0002
$ws.Cells.Item(126, 12).Value = 9000  # L126: 8100 -> 9000
$ws.Cells.Item(126, 13).Value = -2880.0002  # M126: -2865.7142 -> -2880.0002
$ws.Cells.Item(126, 14).Value = -13940  # N126: -13040 -> -13940

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 33493.812  # H132: 31582.53 -> 33493.812
$ws.Cells.Item(132, 9).Value = 1869.1538  # I132: 1753.5333 -> 1869.1538
$ws.Cells.Item(132, 11).Value = 5607.4614  # K132: 5260.5999 -> 5607.4614
$ws.Cells.Item(132, 13).Value = -3077.4614  # M132: -2730.5999 -> -3077.4614
$ws.Cells.Item(135, 8).Value = 35346.727  # H135: 40000 -> 35346.727
$ws.Cells.Item(135, 10).Value = 35346.727  # J135: 40000 -> 35346.727
$ws.Cells.Item(135, 12).Value = 35346.727  # L135: 40000 -> 35346.727
$ws.Cells.Item(135, 14).Value = -45486.727  # N135: -50140 -> -45486.727
$ws.Cells.Item(136, 8).Value = 5578.44  # H136: 5211.2964 -> 5578.44
$ws.Cells.Item(136, 9).Value = 6453.05  # I136: 5709.1304 -> 6453.05
$ws.Cells.Item(136, 10).Value = 2080  # J136: 2348.75 -> 2080
$ws.Cells.Item(136, 11).Value = 19359.15  # K136: 17127.3912 -> 19359.15
$ws.Cells.Item(136, 12).Value = 6240  # L136: 7046.25 -> 6240
$ws.Cells.Item(136, 13).Value = -16809.15  # M136: -14577.3912 -> -16809.15
$ws.Cells.Item(136, 14).Value = -11340  # N136: -12146.25 -> -11340

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(63, 8).Value = 0  # H63: 60000 -> 0
$ws.Cells.Item(63, 10).Value = 0  # J63: 60000 -> 0
$ws.Cells.Item(63, 12).Value = 0  # L63: 60000 -> 0
$ws.Cells.Item(63, 14).ClearContents()  # N63: removed
$ws.Cells.Item(66, 8).Value = 0  # H66: 60000 -> 0
$ws.Cells.Item(66, 10).Value = 0  # J66: 60000 -> 0
$ws.Cells.Item(66, 12).Value = 0  # L66: 180000 -> 0
$ws.Cells.Item(66, 14).ClearContents()  # N66: removed
$ws.Cells.Item(126, 8).Value = 62501740  # H126: 66668420 -> 62501740
$ws.Cells.Item(126, 10).Value = 1833.1666  # J126: 1900 -> 1833.1666
$ws.Cells.Item(126, 12).Value = 5499.4998  # L126: 5700 -> 5499.4998
$ws.Cells.Item(126, 14).Value = -10439.4998  # N126: -10640 -> -10439.4998
$ws.Cells.Item(132, 8).Value = 3490.6765  # H132: 3254.3784 -> 3490.6765
$ws.Cells.Item(132, 9).Value = 4194.1  # I132: 3576.25 -> 4194.1
$ws.Cells.Item(132, 10).Value = 2485.7856  # J132: 2660.1538 -> 2485.7856
$ws.Cells.Item(132, 11).Value = 12582.3  # K132: 10728.75 -> 12582.3
$ws.Cells.Item(132, 12).Value = 7457.3568  # L132: 7980.4614 -> 7457.3568
$ws.Cells.Item(132, 13).Value = -10052.3  # M132: -8198.75 -> -10052.3
$ws.Cells.Item(132, 14).Value = -12517.3568  # N132: -13040.4614 -> -12517.3568

Write-Output "Applied all cell updates."